# feat: add 2022-Q4 data
#
# 1. Insert a brand-new "2022-Q4" worksheet right after "总计" (i.e. before
#    the existing "2022-Q3" tab), populated with the new quarter's fund
#    holdings, formatted like the other quarterly sheets (bold/centered
#    header row + first column, thin borders).
# 2. Update the "总计" (summary) sheet: add a "2022-Q4" row at the top of
#    the data (row 2), push the previously existing rows down by one, and
#    re-append the "2021-Q2" row at the bottom (row 7) so the history stays
#    complete.

$wb = $excel.ActiveWorkbook

function Format-HeaderLike($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1         # xlContinuous
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet before "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

Format-HeaderLike($q4.Range("B1:H1"))
Format-HeaderLike($q4.Range("A1:A3"))

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'002123"
$q4.Range("C2").Value = "北信瑞丰外延增长主题灵活配置混合"
$q4.Range("D2").Value = "'0.20"
$q4.Range("E2").Value = "'93.59"
$q4.Range("F2").Value = "'3.81"
$q4.Range("G2").Value = "'0.0076"
$q4.Range("H2").Value = 7

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'001154"
$q4.Range("C3").Value = "北信瑞丰平安中国主题灵活配置混合"
$q4.Range("D3").Value = "'0.13"
$q4.Range("E3").Value = "'93.42"
$q4.Range("F3").Value = "'3.18"
$q4.Range("G3").Value = "'0.0041"
$q4.Range("H3").Value = 8

$q4.Range("A1:H3").Columns.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet with the new quarter on top, shifting
#    the previously existing rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

Format-HeaderLike($total.Range("A7"))

$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 2
$total.Range("D7").Value = 0.11

$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 1.15

$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 6
$total.Range("D5").Value = 0.59

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 0.27

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.01

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01

$total.Range("A7").Value = 5
